$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 21-33: only the cells that actually changed (per diff) ---
# Row 21
$ws.Range("D21").Value = 44435

# Row 22
$ws.Range("D22").Value = 44327
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = 9500
$ws.Range("P22").Value = 158

# Row 23
$ws.Range("D23").Value = 44313
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = 8500
$ws.Range("P23").Value = 142

# Row 24
$ws.Range("D24").Value = 44350
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = 12000
$ws.Range("L24").Value = 13000
$ws.Range("M24").Value = 12500
$ws.Range("P24").Value = 208

# Row 25
$ws.Range("D25").Value = 44250
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 6500
$ws.Range("P25").Value = 108

# Row 26
$ws.Range("D26").Value = 44294

# Row 27
$ws.Range("D27").Value = 44159
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 8500
$ws.Range("P27").Value = 142

# Row 28
$ws.Range("D28").Value = 44433
$ws.Range("K28").Value = 9000
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 9500
$ws.Range("P28").Value = 158

# Row 29
$ws.Range("D29").Value = 44316
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 9500
$ws.Range("P29").Value = 158

# Row 30
$ws.Range("D30").Value = 44370
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = 12500
$ws.Range("O30").Value = 'Región de Arica y Parinacota'
$ws.Range("P30").Value = 208

# Row 31
$ws.Range("D31").Value = 44230
$ws.Range("K31").Value = 8000
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = 8500
$ws.Range("P31").Value = 142

# Row 32
$ws.Range("D32").Value = 44398
$ws.Range("K32").Value = 13000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 13500
$ws.Range("O32").Value = 'Región Metropolitana'
$ws.Range("P32").Value = 225

# Row 33
$ws.Range("D33").Value = 44355

# --- Append two brand-new rows (34 and 35) at the end ---
# Row 34
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = 'Vega Monumental Concepción'
$ws.Range("C34").Value = 'Bíobío'
$ws.Range("D34").Value = 44223
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112001
$ws.Range("G34").Value = 'Berenjena'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 100
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 8500
$ws.Range("M34").Value = 8250
$ws.Range("N34").Value = '$/caja 60 unidades'
$ws.Range("O34").Value = 'Región de Arica y Parinacota'
$ws.Range("P34").Value = 138
$ws.Range("Q34").Value = 60
$ws.Range("R34").Value = 'Hortaliza'

# Row 35
$ws.Range("A35").Value = 11
$ws.Range("B35").Value = 'Vega Monumental Concepción'
$ws.Range("C35").Value = 'Bíobío'
$ws.Range("D35").Value = 44334
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = 100112001
$ws.Range("G35").Value = 'Berenjena'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 11500
$ws.Range("N35").Value = '$/caja 60 unidades'
$ws.Range("O35").Value = 'Región de Arica y Parinacota'
$ws.Range("P35").Value = 192
$ws.Range("Q35").Value = 60
$ws.Range("R35").Value = 'Hortaliza'

# --- Match the date-column number format used elsewhere in column D ---
$ws.Range("D34").NumberFormat = $ws.Range("D33").NumberFormat
$ws.Range("D35").NumberFormat = $ws.Range("D33").NumberFormat

Write-Output "edit complete"
